$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip trailing whitespace/newlines from the title cells in column A (rows 2-11)
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2.TrimEnd()
}

# Adjust row heights for header row and first data row
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 1274.25
